$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and the Cosmos/InjectiveProtocol row swap)
# Use Text number format so purely-numeric-looking strings are preserved exactly as text
# (matching the original inlineStr cell type) instead of being parsed into numeric values.
$cellUpdates = @{
    'D2' = '43.746.74'
    'E2' = '  +0.40%  '
    'D3' = '2.298.94'
    'E3' = '  +0.45%  '
    'E4' = '  +0.38%  '
    'D5' = '115.50'
    'E5' = '  +21.13%  '
    'D6' = '268.95'
    'E6' = '  +0.40%  '
    'D7' = '0.625'
    'E7' = '  +1.48%  '
    'E8' = '  +0.38%  '
    'D9' = '0.622'
    'E9' = '  +2.27%  '
    'D10' = '48.64'
    'E10' = '  +8.99%  '
    'D11' = '0.0943'
    'E11' = '  +0.52%  '
    'D12' = '8.68'
    'E12' = '  +11.25%  '
    'E13' = '  +2.56%  '
    'D14' = '15.64'
    'E14' = '  +3.08%  '
    'D15' = '2.646.57'
    'E15' = '  +0.53%  '
    'D16' = '0.859'
    'E16' = '  +0.92%  '
    'D17' = '2.301.43'
    'E17' = '  +0.47%  '
    'D18' = '43.641.60'
    'E18' = '  +0.21%  '
    'D19' = '0.0000110'
    'E19' = '  +2.53%  '
    'D20' = '6.58'
    'E20' = '  +6.08%  '
    'D21' = '72.77'
    'E21' = '  +0.43%  '
    'D22' = '2.55'
    'E22' = '  +3.81%  '
    'D23' = '234.14'
    'E23' = '  -0.42%  '
    'D24' = '9.65'
    'E24' = '  +6.75%  '
    'E25' = '  +14.86%  '
    'E26' = '  -0.06%  '
    'B27' = 'InjectiveProtocol'
    'C27' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D27' = '44.60'
    'E27' = '  +11.44%  '
    'B28' = 'Cosmos'
    'C28' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D28' = '11.51'
    'E28' = '  +2.56%  '
    'E29' = '  -1.42%  '
    'E30' = '  -0.29%  '
    'D31' = '177.45'
    'E31' = '  +1.30%  '
    'D32' = '0.0939'
    'E32' = '  +6.45%  '
    'D33' = '21.79'
    'E33' = '  -0.88%  '
    'D34' = '5.55'
    'E34' = '  +3.62%  '
    'D35' = '0.126'
    'E35' = '  +0.92%  '
    'E36' = '  +8.24%  '
    'E37' = '  +1.84%  '
    'D38' = '3.93'
    'E38' = '  +18.75%  '
    'D39' = '0.0357'
    'E39' = '  +0.07%  '
    'D40' = '75.28'
    'E40' = '  +16.58%  '
    'D41' = '0.243'
    'E41' = '  +3.11%  '
    'D42' = '2.41'
    'E42' = '  +3.05%  '
    'D43' = '13.31'
    'E43' = '  +10.69%  '
    'D44' = '1.43'
    'E44' = '  +6.79%  '
    'E45' = '  +0.30%  '
    'D46' = '5.96'
    'E46' = '  +13.94%  '
    'D47' = '8.82'
    'E47' = '  +0.03%  '
    'E48' = '  -1.18%  '
    'D49' = '101.86'
    'E49' = '  +3.77%  '
    'E50' = '  +4.43%  '
    'D51' = '0.454'
    'E51' = '  +5.82%  '
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
}
